# Apply cryptocurrency price/volume updates per the commit diff.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextCell($cell, $text) {
    # Force the cell to hold a literal text string (matches the source
    # workbook, where these are inlineStr cells, not numbers) while
    # restoring the default "Normal" style afterwards so no stray
    # cell-level style index is left behind.
    $cell.NumberFormat = "@"
    $cell.Value = $text
    $cell.Style = "Normal"
}

Set-TextCell $ws.Range("D2") "69.388.25"
$ws.Range("E2").Value = "  +1.77%  "

Set-TextCell $ws.Range("D3") "3.944.41"
$ws.Range("E3").Value = "  +0.66%  "

$ws.Range("E4").Value = "  -0.12%  "

Set-TextCell $ws.Range("D5") "493.19"
$ws.Range("E5").Value = "  +1.02%  "

Set-TextCell $ws.Range("D6") "147.76"
$ws.Range("E6").Value = "  +0.86%  "

Set-TextCell $ws.Range("D7") "0.624"
$ws.Range("E7").Value = "  -0.44%  "

Set-TextCell $ws.Range("D9") "0.735"
$ws.Range("E9").Value = "  +0.61%  "

$ws.Range("E10").Value = "  +3.58%  "

Set-TextCell $ws.Range("D11") "0.0000350"
$ws.Range("E11").Value = "  -1.63%  "

Set-TextCell $ws.Range("D12") "43.32"
$ws.Range("E12").Value = "  +1.41%  "

Set-TextCell $ws.Range("D13") "10.44"
$ws.Range("E13").Value = "  -1.75%  "

Set-TextCell $ws.Range("D14") "4.576.83"
$ws.Range("E14").Value = "  +0.70%  "

Set-TextCell $ws.Range("D15") "3.974.22"
$ws.Range("E15").Value = "  +1.35%  "

Set-TextCell $ws.Range("D16") "14.35"
$ws.Range("E16").Value = "  -2.82%  "

$ws.Range("E17").Value = "  -0.75%  "

Set-TextCell $ws.Range("D18") "19.86"
$ws.Range("E18").Value = "  -0.68%  "

$ws.Range("E19").Value = "  +2.69%  "

Set-TextCell $ws.Range("D20") "69.460.49"
$ws.Range("E20").Value = "  +1.67%  "

Set-TextCell $ws.Range("D21") "439.93"
$ws.Range("E21").Value = "  -0.74%  "

Set-TextCell $ws.Range("D22") "3.46"
$ws.Range("E22").Value = "  +2.80%  "

Set-TextCell $ws.Range("D23") "14.53"
$ws.Range("E23").Value = "  -1.36%  "

Set-TextCell $ws.Range("D24") "88.97"
$ws.Range("E24").Value = "  +0.63%  "

Set-TextCell $ws.Range("D25") "12.04"
$ws.Range("E25").Value = "  +8.54%  "

Set-TextCell $ws.Range("D26") "3.78"
$ws.Range("E26").Value = "  +3.82%  "

Set-TextCell $ws.Range("D27") "11.16"
$ws.Range("E27").Value = "  -4.60%  "

Set-TextCell $ws.Range("D28") "37.20"
$ws.Range("E28").Value = "  -4.22%  "

Set-TextCell $ws.Range("D29") "5.63"
$ws.Range("E29").Value = "  -4.33%  "

Set-TextCell $ws.Range("D30") "707.63"
$ws.Range("E30").Value = "  -0.28%  "

$ws.Range("E31").Value = "  +0.65%  "

Set-TextCell $ws.Range("D32") "13.33"
$ws.Range("E32").Value = "  -1.10%  "

Set-TextCell $ws.Range("D33") "2.89"
$ws.Range("E33").Value = "  +0.85%  "

Set-TextCell $ws.Range("D34") "0.470"
$ws.Range("E34").Value = "  +23.78%  "

$ws.Range("E35").Value = "  +1.73%  "

Set-TextCell $ws.Range("D38") "40.87"
$ws.Range("E38").Value = "  -1.21%  "

Set-TextCell $ws.Range("D39") "0.151"
$ws.Range("E39").Value = "  +0.98%  "

Set-TextCell $ws.Range("D40") "0.998"
$ws.Range("E40").Value = "  -0.10%  "

$ws.Range("E41").Value = "  +0.09%  "

Set-TextCell $ws.Range("D42") "0.0492"
$ws.Range("E42").Value = "  +2.59%  "

Set-TextCell $ws.Range("D43") "2.96"
$ws.Range("E43").Value = "  +3.62%  "

Set-TextCell $ws.Range("D44") "3.06"
$ws.Range("E44").Value = "  -1.77%  "

$ws.Range("E45").Value = "  +2.27%  "

$ws.Range("E46").Value = "  +0.55%  "

$ws.Range("E47").Value = "  +6.64%  "

Set-TextCell $ws.Range("D48") "0.0₆0359"
$ws.Range("E48").Value = "  +4.42%  "

Set-TextCell $ws.Range("D49") "3.04"
$ws.Range("E49").Value = "  +7.56%  "

$ws.Range("E50").Value = "  -0.43%  "

Set-TextCell $ws.Range("D51") "143.83"
$ws.Range("E51").Value = "  -0.75%  "

# Rows 36/37: OKB and NEARProtocol swap rank positions; new price/volume values.
$ws.Range("B36").Value = "NEARProtocol"
$ws.Range("C36").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
Set-TextCell $ws.Range("D36") "6.08"
$ws.Range("E36").Value = "  +4.75%  "

$ws.Range("B37").Value = "OKB"
$ws.Range("C37").Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
Set-TextCell $ws.Range("D37") "61.54"
$ws.Range("E37").Value = "  +4.25%  "
